$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 815756
$ws.Range("R2").Value = 7420847
$ws.Range("Z2").Value = $null
$ws.Range("AB2").Value = $null
